$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string cell updates (Coin names, Links, prices with thousand-dot separators, volume %% strings)
$ws.Cells.Item(2, 4).Value = "51.629.40"
$ws.Cells.Item(2, 5).Value = "  +1.41%  "
$ws.Cells.Item(3, 4).Value = "2.993.36"
$ws.Cells.Item(3, 5).Value = "  +1.91%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 5).Value = "  +1.78%  "
$ws.Cells.Item(6, 5).Value = "  +3.87%  "
$ws.Cells.Item(7, 5).Value = "  +1.57%  "
$ws.Cells.Item(8, 5).Value = "  -0.01%  "
$ws.Cells.Item(9, 5).Value = "  +1.93%  "
$ws.Cells.Item(10, 5).Value = "  +1.57%  "
$ws.Cells.Item(11, 5).Value = "  -0.59%  "
$ws.Cells.Item(12, 5).Value = "  +1.66%  "
$ws.Cells.Item(13, 4).Value = "3.464.36"
$ws.Cells.Item(13, 5).Value = "  +2.19%  "
$ws.Cells.Item(14, 5).Value = "  +3.81%  "
$ws.Cells.Item(15, 5).Value = "  +2.63%  "
$ws.Cells.Item(16, 4).Value = "2.980.00"
$ws.Cells.Item(16, 5).Value = "  +1.60%  "
$ws.Cells.Item(17, 5).Value = "  -0.52%  "
$ws.Cells.Item(18, 5).Value = "  +0.50%  "
$ws.Cells.Item(19, 4).Value = "51.634.78"
$ws.Cells.Item(19, 5).Value = "  +1.49%  "
$ws.Cells.Item(20, 5).Value = "  +1.38%  "
$ws.Cells.Item(21, 5).Value = "  +1.33%  "
$ws.Cells.Item(22, 5).Value = "  +1.04%  "
$ws.Cells.Item(23, 5).Value = "  +2.15%  "
$ws.Cells.Item(24, 5).Value = "  +0.92%  "
$ws.Cells.Item(25, 5).Value = "  +2.11%  "
$ws.Cells.Item(26, 5).Value = "  -1.45%  "
$ws.Cells.Item(27, 2).Value = "Kaspa"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(27, 5).Value = "  +5.31%  "
$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(28, 5).Value = "  -1.57%  "
$ws.Cells.Item(29, 5).Value = "  +0.05%  "
$ws.Cells.Item(30, 5).Value = "  +2.06%  "
$ws.Cells.Item(31, 5).Value = "  +1.10%  "
$ws.Cells.Item(32, 5).Value = "  +4.38%  "
$ws.Cells.Item(33, 5).Value = "  +4.29%  "
$ws.Cells.Item(34, 5).Value = "  +1.04%  "
$ws.Cells.Item(35, 5).Value = "  +0.53%  "
$ws.Cells.Item(36, 5).Value = "  +1.68%  "
$ws.Cells.Item(37, 5).Value = "  +0.02%  "
$ws.Cells.Item(38, 5).Value = "  +5.98%  "
$ws.Cells.Item(39, 5).Value = "  +3.14%  "
$ws.Cells.Item(40, 5).Value = "  +5.47%  "
$ws.Cells.Item(41, 2).Value = "Stellar"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(41, 5).Value = "  +1.49%  "
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 5).Value = "  +2.32%  "
$ws.Cells.Item(43, 2).Value = "Monero"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(43, 5).Value = "  +6.20%  "
$ws.Cells.Item(44, 2).Value = "NEARProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(44, 5).Value = "  +13.64%  "
$ws.Cells.Item(45, 5).Value = "  +0.98%  "
$ws.Cells.Item(46, 5).Value = "  -0.39%  "
$ws.Cells.Item(47, 5).Value = "  +0.35%  "
$ws.Cells.Item(48, 5).Value = "  +0.75%  "
$ws.Cells.Item(49, 4).Value = "2.039.61"
$ws.Cells.Item(49, 5).Value = "  +2.60%  "
$ws.Cells.Item(50, 4).Value = "3.289.09"
$ws.Cells.Item(50, 5).Value = "  +1.75%  "
$ws.Cells.Item(51, 5).Value = "  +1.19%  "

# Numeric-looking price values that must remain stored as text (force text format, then clear the
# temporary formatting so the cell keeps its original default style, matching source formatting).
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "381.87"
$cell.ClearFormats()
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "104.48"
$cell.ClearFormats()
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "36.75"
$cell.ClearFormats()
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.86"
$cell.ClearFormats()
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.55"
$cell.ClearFormats()
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.18"
$cell.ClearFormats()
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.62"
$cell.ClearFormats()
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.52"
$cell.ClearFormats()
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "267.96"
$cell.ClearFormats()
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.171"
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.27"
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.16"
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "34.68"
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.43"
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.06"
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0448"
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.31"
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.117"
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.85"
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "127.43"
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.85"
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.37"
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0331"
$cell.ClearFormats()
